# "complete monthly and re-run daily"
# Re-running the daily/monthly productivity pipeline changed the
# arrive_thres / count / prod_*_ratio figures for employee sp99004
# (陳衣玲) across the four output sheets.

$wb = $excel.ActiveWorkbook

# --- team_df (sheet1): two individual shift rows for sp99004 ---
$wsTeamDf = $wb.Worksheets.Item("team_df")

# Row 6 -> arrive_thres 2->3, count 10->16, prod_hour_ratio 0.2->0.1875
$wsTeamDf.Range("S6").Value = 3
$wsTeamDf.Range("T6").Value = 16
$wsTeamDf.Range("U6").Value = 0.1875

# Row 12 -> count 6->10, prod_hour_ratio 0.3333333333333333->0.2
$wsTeamDf.Range("T12").Value = 10
$wsTeamDf.Range("U12").Value = 0.2

# --- team_df_day (sheet2): daily roll-up for sp99004 ---
$wsTeamDfDay = $wb.Worksheets.Item("team_df_day")

# Row 2 -> arrive_thres 4->5, count 16->26, prod_day_ratio 0.25->0.1923076923076923
$wsTeamDfDay.Range("F2").Value = 5
$wsTeamDfDay.Range("G2").Value = 26
$wsTeamDfDay.Range("H2").Value = 0.1923076923076923

# --- productivity_tl (sheet3): TL productivity score for sp99004 ---
$wsProdTl = $wb.Worksheets.Item("productivity_tl")
$wsProdTl.Range("D2").Value = 0.1923076923076923

# --- productivity_team_function (sheet4): same score for sp99004 ---
$wsProdTeamFunction = $wb.Worksheets.Item("productivity_team_function")
$wsProdTeamFunction.Range("D2").Value = 0.1923076923076923
